$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Thbs1"
$ws.Cells.Item(2,3).Value = "Tnfrsf11b"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3.0
$ws.Cells.Item(2,6).Value = 1.0
$ws.Cells.Item(2,7).Value = 31.29063333333333
$ws.Cells.Item(2,8).Value = 93.8719
$ws.Cells.Item(2,9).Value = 0.02026792284095206
$ws.Cells.Item(2,10).Value = 0.02026792284095206
$ws.Cells.Item(2,11).Value = 1.0
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.08241233333333334
$ws.Cells.Item(2,14).Value = 0.247237
$ws.Cells.Item(2,15).Value = 0.04727005612861496
$ws.Cells.Item(2,16).Value = 0.04727005612861496
$ws.Cells.Item(2,17).Value = 2.578734104477778
$ws.Cells.Item(2,18).Value = 23.2086069403
$ws.Cells.Item(2,19).Value = 0.0009580658503022409
$ws.Cells.Item(2,20).Value = 0.0009580658503022409

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Thbs1"
$ws.Cells.Item(3,3).Value = "Tnfrsf11b"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3.0
$ws.Cells.Item(3,6).Value = 1.0
$ws.Cells.Item(3,7).Value = 31.29063333333333
$ws.Cells.Item(3,8).Value = 93.8719
$ws.Cells.Item(3,9).Value = 0.02026792284095206
$ws.Cells.Item(3,10).Value = 0.02026792284095206
$ws.Cells.Item(3,11).Value = 3.0
$ws.Cells.Item(3,12).Value = 1.0
$ws.Cells.Item(3,13).Value = 1.661024
$ws.Cells.Item(3,14).Value = 4.983072
$ws.Cells.Item(3,15).Value = 0.9527299438713851
$ws.Cells.Item(3,16).Value = 0.952729943871385
$ws.Cells.Item(3,17).Value = 51.97449294186666
$ws.Cells.Item(3,18).Value = 467.7704364768
$ws.Cells.Item(3,19).Value = 0.01930985699064982
$ws.Cells.Item(3,20).Value = 0.01930985699064981

# Row 4
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Thbs1"
$ws.Cells.Item(4,3).Value = "Tnfrsf11b"
$ws.Cells.Item(4,4).Value = "ECs"
$ws.Cells.Item(4,5).Value = 3.0
$ws.Cells.Item(4,6).Value = 1.0
$ws.Cells.Item(4,7).Value = 170.232249
$ws.Cells.Item(4,8).Value = 510.696747
$ws.Cells.Item(4,9).Value = 0.1102647572204378
$ws.Cells.Item(4,10).Value = 0.1102647572204378
$ws.Cells.Item(4,11).Value = 1.0
$ws.Cells.Item(4,12).Value = 0.3333333333333333
$ws.Cells.Item(4,13).Value = 0.08241233333333334
$ws.Cells.Item(4,14).Value = 0.247237
$ws.Cells.Item(4,15).Value = 0.04727005612861496
$ws.Cells.Item(4,16).Value = 0.04727005612861496
$ws.Cells.Item(4,17).Value = 14.029236848671
$ws.Cells.Item(4,18).Value = 126.263131638039
$ws.Cells.Item(4,19).Value = 0.005212221262818196
$ws.Cells.Item(4,20).Value = 0.005212221262818196

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Thbs1"
$ws.Cells.Item(5,3).Value = "Tnfrsf11b"
$ws.Cells.Item(5,4).Value = "FAPs"
$ws.Cells.Item(5,5).Value = 3.0
$ws.Cells.Item(5,6).Value = 1.0
$ws.Cells.Item(5,7).Value = 170.232249
$ws.Cells.Item(5,8).Value = 510.696747
$ws.Cells.Item(5,9).Value = 0.1102647572204378
$ws.Cells.Item(5,10).Value = 0.1102647572204378
$ws.Cells.Item(5,11).Value = 3.0
$ws.Cells.Item(5,12).Value = 1.0
$ws.Cells.Item(5,13).Value = 1.661024
$ws.Cells.Item(5,14).Value = 4.983072
$ws.Cells.Item(5,15).Value = 0.9527299438713851
$ws.Cells.Item(5,16).Value = 0.952729943871385
$ws.Cells.Item(5,17).Value = 282.759851162976
$ws.Cells.Item(5,18).Value = 2544.838660466784
$ws.Cells.Item(5,19).Value = 0.1050525359576196
$ws.Cells.Item(5,20).Value = 0.1050525359576196

# Row 6
$ws.Cells.Item(6,1).Value = "M1"
$ws.Cells.Item(6,2).Value = "Thbs1"
$ws.Cells.Item(6,3).Value = "Tnfrsf11b"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3.0
$ws.Cells.Item(6,6).Value = 1.0
$ws.Cells.Item(6,7).Value = 603.9765116666666
$ws.Cells.Item(6,8).Value = 1811.929535
$ws.Cells.Item(6,9).Value = 0.391214495590503
$ws.Cells.Item(6,10).Value = 0.391214495590503
$ws.Cells.Item(6,11).Value = 1.0
$ws.Cells.Item(6,12).Value = 0.3333333333333333
$ws.Cells.Item(6,13).Value = 0.08241233333333334
$ws.Cells.Item(6,14).Value = 0.247237
$ws.Cells.Item(6,15).Value = 0.04727005612861496
$ws.Cells.Item(6,16).Value = 0.04727005612861496
$ws.Cells.Item(6,17).Value = 49.77511360497722
$ws.Cells.Item(6,18).Value = 447.976022444795
$ws.Cells.Item(6,19).Value = 0.01849273116489087
$ws.Cells.Item(6,20).Value = 0.01849273116489087

# Row 7
$ws.Cells.Item(7,1).Value = "M1"
$ws.Cells.Item(7,2).Value = "Thbs1"
$ws.Cells.Item(7,3).Value = "Tnfrsf11b"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3.0
$ws.Cells.Item(7,6).Value = 1.0
$ws.Cells.Item(7,7).Value = 603.9765116666666
$ws.Cells.Item(7,8).Value = 1811.929535
$ws.Cells.Item(7,9).Value = 0.391214495590503
$ws.Cells.Item(7,10).Value = 0.391214495590503
$ws.Cells.Item(7,11).Value = 3.0
$ws.Cells.Item(7,12).Value = 1.0
$ws.Cells.Item(7,13).Value = 1.661024
$ws.Cells.Item(7,14).Value = 4.983072
$ws.Cells.Item(7,15).Value = 0.9527299438713851
$ws.Cells.Item(7,16).Value = 0.952729943871385
$ws.Cells.Item(7,17).Value = 1003.219481314613
$ws.Cells.Item(7,18).Value = 9028.97533183152
$ws.Cells.Item(7,19).Value = 0.3727217644256121
$ws.Cells.Item(7,20).Value = 0.3727217644256121

# Row 8
$ws.Cells.Item(8,1).Value = "M2"
$ws.Cells.Item(8,2).Value = "Thbs1"
$ws.Cells.Item(8,3).Value = "Tnfrsf11b"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3.0
$ws.Cells.Item(8,6).Value = 1.0
$ws.Cells.Item(8,7).Value = 694.886846
$ws.Cells.Item(8,8).Value = 2084.660538
$ws.Cells.Item(8,9).Value = 0.4500999653118942
$ws.Cells.Item(8,10).Value = 0.4500999653118942
$ws.Cells.Item(8,11).Value = 1.0
$ws.Cells.Item(8,12).Value = 0.3333333333333333
$ws.Cells.Item(8,13).Value = 0.08241233333333334
$ws.Cells.Item(8,14).Value = 0.247237
$ws.Cells.Item(8,15).Value = 0.04727005612861496
$ws.Cells.Item(8,16).Value = 0.04727005612861496
$ws.Cells.Item(8,17).Value = 57.26724638150067
$ws.Cells.Item(8,18).Value = 515.4052174335061
$ws.Cells.Item(8,19).Value = 0.02127625062378089
$ws.Cells.Item(8,20).Value = 0.02127625062378089

# Row 9
$ws.Cells.Item(9,1).Value = "M2"
$ws.Cells.Item(9,2).Value = "Thbs1"
$ws.Cells.Item(9,3).Value = "Tnfrsf11b"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3.0
$ws.Cells.Item(9,6).Value = 1.0
$ws.Cells.Item(9,7).Value = 694.886846
$ws.Cells.Item(9,8).Value = 2084.660538
$ws.Cells.Item(9,9).Value = 0.4500999653118942
$ws.Cells.Item(9,10).Value = 0.4500999653118942
$ws.Cells.Item(9,11).Value = 3.0
$ws.Cells.Item(9,12).Value = 1.0
$ws.Cells.Item(9,13).Value = 1.661024
$ws.Cells.Item(9,14).Value = 4.983072
$ws.Cells.Item(9,15).Value = 0.9527299438713851
$ws.Cells.Item(9,16).Value = 0.952729943871385
$ws.Cells.Item(9,17).Value = 1154.223728490304
$ws.Cells.Item(9,18).Value = 10388.01355641274
$ws.Cells.Item(9,19).Value = 0.4288237146881134
$ws.Cells.Item(9,20).Value = 0.4288237146881133

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Thbs1"
$ws.Cells.Item(10,3).Value = "Tnfrsf11b"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 3.0
$ws.Cells.Item(10,6).Value = 1.0
$ws.Cells.Item(10,7).Value = 43.46379233333334
$ws.Cells.Item(10,8).Value = 130.391377
$ws.Cells.Item(10,9).Value = 0.02815285903621308
$ws.Cells.Item(10,10).Value = 0.02815285903621308
$ws.Cells.Item(10,11).Value = 1.0
$ws.Cells.Item(10,12).Value = 0.3333333333333333
$ws.Cells.Item(10,13).Value = 0.08241233333333334
$ws.Cells.Item(10,14).Value = 0.247237
$ws.Cells.Item(10,15).Value = 0.04727005612861496
$ws.Cells.Item(10,16).Value = 0.04727005612861496
$ws.Cells.Item(10,17).Value = 3.581952541705445
$ws.Cells.Item(10,18).Value = 32.237572875349
$ws.Cells.Item(10,19).Value = 0.001330787226822778
$ws.Cells.Item(10,20).Value = 0.001330787226822777

# Row 11
$ws.Cells.Item(11,1).Value = "sCs"
$ws.Cells.Item(11,2).Value = "Thbs1"
$ws.Cells.Item(11,3).Value = "Tnfrsf11b"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 3.0
$ws.Cells.Item(11,6).Value = 1.0
$ws.Cells.Item(11,7).Value = 43.46379233333334
$ws.Cells.Item(11,8).Value = 130.391377
$ws.Cells.Item(11,9).Value = 0.02815285903621308
$ws.Cells.Item(11,10).Value = 0.02815285903621308
$ws.Cells.Item(11,11).Value = 3.0
$ws.Cells.Item(11,12).Value = 1.0
$ws.Cells.Item(11,13).Value = 1.661024
$ws.Cells.Item(11,14).Value = 4.983072
$ws.Cells.Item(11,15).Value = 0.9527299438713851
$ws.Cells.Item(11,16).Value = 0.952729943871385
$ws.Cells.Item(11,17).Value = 72.19440219668267
$ws.Cells.Item(11,18).Value = 649.749619770144
$ws.Cells.Item(11,19).Value = 0.02682207180939031
$ws.Cells.Item(11,20).Value = 0.0268220718093903
